$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new company record (CSS / IRPJ lookup result).
$ws.Range("A2").Value = "2W COBERTURAS LTDA"

# CNPJ must stay textual (it has a meaningful leading/trailing digit pattern,
# not a numeric quantity), so force the cell to Text format before writing it.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "40567022000152"

# The previous second result row is no longer part of the saved output -
# remove it entirely so the used range shrinks back down to A1:B2.
$ws.Rows.Item(3).Delete()

# Leave the selection parked back on the header cell, as it was when saved.
$ws.Range("A1").Select() | Out-Null
